$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 1, shifting all existing rows down by one.
$ws.Rows.Item(1).Insert()

# Set the new cell's value.
$ws.Range("A1").Value = "GO_Biological_Process_2023"

# Match the final selection state recorded in the saved file.
$ws.Range("B8").Select()
